$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptos snapshot: Price (D) and Volume(1h) (E) columns.
# D-column values that parse as plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matches source data like
# "300.64" rather than letting them become floating point numbers), then
# ClearFormats() strips the quote-prefix styling that introduces so the
# cell keeps its original (default) style.

$ws.Range("D2").Value = "46.032.32"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.340.19"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'300.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "'98.59"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("E9").Value = "  -5.83%  "
$ws.Range("D10").Value = "'34.58"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").Value = "'7.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.04%  "
$ws.Range("E13").Value = "  -1.43%  "
$ws.Range("D14").Value = "2.695.31"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "2.338.23"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  -3.43%  "
$ws.Range("E17").Value = "  -4.78%  "
$ws.Range("D18").Value = "46.038.59"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "'12.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -8.22%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'5.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").Value = "'66.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").Value = "'244.69"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("E24").Value = "  -5.96%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -5.39%  "
$ws.Range("D27").Value = "'39.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.72%  "
$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").Value = "'21.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "'3.60"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +15.10%  "
$ws.Range("D32").Value = "'2.81"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.94%  "
$ws.Range("D33").Value = "'5.46"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("D34").Value = "'144.15"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.62%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("E37").Value = "  -3.02%  "
$ws.Range("D38").Value = "'1.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "'14.93"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").Value = "'3.86"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("E42").Value = "  -6.01%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "1.838.22"
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").Value = "'90.24"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").Value = "'1.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.61%  "
$ws.Range("E47").Value = "  -5.26%  "
$ws.Range("D48").Value = "'69.09"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.88%  "
$ws.Range("D49").Value = "2.567.86"
$ws.Range("D50").Value = "'96.39"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").Value = "'8.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.48%  "
